$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells are stored as text (avoid numeric/date auto-conversion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.493.96'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.853.20'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '240.63'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '0.6325'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '0.07584'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').Value = '0.2966'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').Value = '24.74'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '0.07703'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '1.851.14'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '5.007'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '0.6848'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').Value = '0.00001027'
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('D16').Value = '83.42'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '2.110.70'
$ws.Range('E17').Value = '  -4.12%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '6.160'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '29.514.64'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '229.39'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '12.52'
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '0.9998'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '7.545'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = '0.9999'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '156.75'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '0.1402'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '8.394'
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '17.71'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.472'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.271'
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.05685'
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.133'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.037'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '1.849'
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.161'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7172'
$ws.Range('E36').Value = '  -0.97%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.592'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.250.20'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01809'
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.780'
$ws.Range('E40').Value = '  -0.90%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '0.9117'
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '6.205'
$ws.Range('E42').Value = '  +1.07%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '101.83'
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '66.59'
$ws.Range('E45').Value = '  -2.11%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('E46').Value = '  -5.17%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '7.092'
$ws.Range('E47').Value = '  -3.92%  '
$ws.Range('D48').Value = '0.4032'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '9.126'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '1.704'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.1125'
$ws.Range('E51').Value = '  +0.58%  '
